$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new headers for the team's win/loss/tie record
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the header formatting used by the rest of row 1 (bold, bordered, centered)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's record (104-58-0) for every player row
for ($r = 2; $r -le 54; $r++) {
    $ws.Cells.Item($r, 30).Value = 104   # column AD - Wins
    $ws.Cells.Item($r, 31).Value = 58    # column AE - Losses
    $ws.Cells.Item($r, 32).Value = 0     # column AF - Ties
}
